$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "43.920.18"
$ws.Range("E2").Value = "  +0.43%  "

# Row 3
$ws.Range("D3").Value = "2.238.10"
$ws.Range("E3").Value = "  -1.31%  "

# Row 4
$ws.Range("E4").Value = "  +0.26%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.06%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.56"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.55%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.568"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.23%  "

# Row 8
$ws.Range("E8").Value = "  +0.08%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.534"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -6.42%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.21"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.36%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0821"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.32%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.36"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -6.24%  "

# Row 13
$ws.Range("E13").Value = "  -3.08%  "

# Row 14
$ws.Range("D14").Value = "2.579.15"
$ws.Range("E14").Value = "  -1.18%  "

# Row 15
$ws.Range("D15").Value = "2.247.48"
$ws.Range("E15").Value = "  -0.65%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.839"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.08%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.09"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.06%  "

# Row 18
$ws.Range("D18").Value = "43.745.21"
$ws.Range("E18").Value = "  +0.11%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.03"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -8.31%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0963"
$ws.Range("E20").Value = "  -3.03%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.98%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "64.93"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.70%  "

# Row 23
$ws.Range("E23").Value = "  -5.52%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "233.41"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.21%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.04"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -8.28%  "

# Row 26
$ws.Range("E26").Value = "  +0.34%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.21"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.03%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.17"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.39%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "36.82"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.74%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.96"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -7.81%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "158.27"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.95%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.95"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.53%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0830"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.94%  "

# Row 34
$ws.Range("E34").Value = "  -1.16%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.18"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.48%  "

# Row 36
$ws.Range("E36").Value = "  +3.11%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.90"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.74%  "

# Row 38
$ws.Range("E38").Value = "  -3.54%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "15.88"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.79%  "

# Row 40
$ws.Range("E40").Value = "  -7.34%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.06"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -9.60%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0309"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.48%  "

# Row 43
$ws.Range("E43").Value = "  +0.07%  "

# Row 44
$ws.Range("D44").Value = "1.729.05"
$ws.Range("E44").Value = "  -5.38%  "

# Row 45
$ws.Range("E45").Value = "  -6.57%  "

# Row 46
$ws.Range("B46").Value = "BitcoinSV"
$ws.Range("C46").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "80.72"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.92%  "

# Row 47
$ws.Range("B47").Value = "ordi"
$ws.Range("C47").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "73.37"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.59%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.09"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.53%  "

# Row 49
$ws.Range("E49").Value = "  +1.36%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "101.65"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.11%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "56.44"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.36%  "
